$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename "men/" -> "men1/" for the existing data rows (2-31, cols A & B)
for ($r = 2; $r -le 31; $r++) {
    $a = $ws.Cells.Item($r, 1).Value2
    $b = $ws.Cells.Item($r, 2).Value2
    $ws.Cells.Item($r, 1).Value = $a.Replace("men/", "men1/")
    $ws.Cells.Item($r, 2).Value = $b.Replace("men/", "men1/")
}

# Completed experiment 2 readings: append 4 more rows (32-35)
$ws.Cells.Item(32, 1).Value = "men1/31-original.jpg"
$ws.Cells.Item(32, 2).Value = "men1/31-scramble.jpg"
$ws.Cells.Item(33, 1).Value = "men1/32-original.jpg"
$ws.Cells.Item(33, 2).Value = "men1/32-scramble.jpg"
$ws.Cells.Item(34, 1).Value = "men1/33-original.jpg"
$ws.Cells.Item(34, 2).Value = "men1/33-scramble.jpg"
$ws.Cells.Item(35, 1).Value = "men1/34-original.jpg"
$ws.Cells.Item(35, 2).Value = "men1/34-scramble.jpg"

# Update selection to match the author's final cursor state
[void]$ws.Range("A2:A35").Select()

# Autofit the two data columns to their content (matches bestFit width change)
$ws.Columns("A:B").AutoFit()
